$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1393.25
$ws.Range("I80").Value = 1516.6666
$ws.Range("K80").Value = 4549.9998
$ws.Range("M80").Value = -3551.9998
$ws.Range("H83").Value = 1393.25
$ws.Range("I83").Value = 1516.6666
$ws.Range("K83").Value = 13649.9994
$ws.Range("M83").Value = -8657.999400000001
$ws.Range("H100").Value = 1131.75
$ws.Range("I100").Value = 1131.75
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1131.75
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -590.75
$ws.Range("H113").Value = 6496.6665
$ws.Range("I113").Value = 2990
$ws.Range("J113").Value = 8250
$ws.Range("K113").Value = 2990
$ws.Range("L113").Value = 8250
$ws.Range("M113").Value = 264
$ws.Range("N113").Value = -14758
$ws.Range("H138").Value = 3353
$ws.Range("I138").Value = 3149.5
$ws.Range("K138").Value = 9448.5
$ws.Range("M138").Value = -4308.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3637.2917
$ws.Range("I74").Value = 2957.4707
$ws.Range("J74").Value = 5288.2856
$ws.Range("K74").Value = 2957.4707
$ws.Range("L74").Value = 5288.2856
$ws.Range("M74").Value = -2083.4707
$ws.Range("N74").Value = -7036.2856
$ws.Range("H77").Value = 3637.2917
$ws.Range("I77").Value = 2957.4707
$ws.Range("J77").Value = 5288.2856
$ws.Range("K77").Value = 14787.3535
$ws.Range("L77").Value = 26441.428
$ws.Range("M77").Value = -10419.3535
$ws.Range("N77").Value = -35177.428
$ws.Range("H92").Value = 36850
$ws.Range("J92").Value = 36850
$ws.Range("L92").Value = 36850
$ws.Range("N92").Value = -41842
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H122").Value = 3468
$ws.Range("I122").Value = 3210.2856
$ws.Range("J122").Value = 4370
$ws.Range("K122").Value = 9630.856800000001
$ws.Range("L122").Value = 13110
$ws.Range("M122").Value = -7180.856800000001
$ws.Range("N122").Value = -18010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2628.9546
$ws.Range("I105").Value = 2448
$ws.Range("J105").Value = 3443.25
$ws.Range("K105").Value = 2448
$ws.Range("L105").Value = 3443.25
$ws.Range("M105").Value = -701
$ws.Range("N105").Value = -6937.25
$ws.Range("H134").Value = 1896.2727
$ws.Range("I134").Value = 1929.3334
$ws.Range("J134").Value = 1747.5
$ws.Range("K134").Value = 5788.0002
$ws.Range("L134").Value = 5242.5
$ws.Range("M134").Value = -3253.0002
$ws.Range("N134").Value = -10312.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4275.5
$ws.Range("I94").Value = 2236.2856
$ws.Range("J94").Value = 7130.4
$ws.Range("K94").Value = 2236.2856
$ws.Range("L94").Value = 7130.4
$ws.Range("M94").Value = -1785.2856
$ws.Range("N94").Value = -8032.4
$ws.Range("H99").Value = 2300
$ws.Range("I99").Value = 1950
$ws.Range("J99").Value = 2387.5
$ws.Range("K99").Value = 1950
$ws.Range("L99").Value = 2387.5
$ws.Range("M99").Value = -452
$ws.Range("N99").Value = -5383.5
$ws.Range("H126").Value = 2300
$ws.Range("I126").Value = 1950
$ws.Range("J126").Value = 2387.5
$ws.Range("K126").Value = 5850
$ws.Range("L126").Value = 7162.5
$ws.Range("M126").Value = -3380
$ws.Range("N126").Value = -12102.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 214.66667
$ws.Range("I12").Value = 189.33333
$ws.Range("J12").Value = 227.33333
$ws.Range("K12").Value = 567.99999
$ws.Range("L12").Value = 681.99999
$ws.Range("M12").Value = -394.99999
$ws.Range("N12").Value = -1027.99999
$ws.Range("H18").Value = 2965
$ws.Range("I18").Value = 1001
$ws.Range("K18").Value = 3003
$ws.Range("M18").Value = -2834
$ws.Range("H39").Value = 6334.5386
$ws.Range("J39").Value = 7429
$ws.Range("L39").Value = 22287
$ws.Range("N39").Value = -22875
$ws.Range("H52").Value = 500
$ws.Range("J52").Value = 500
$ws.Range("L52").Value = 1500
$ws.Range("N52").Value = -2032
$ws.Range("H57").Value = 1299.125
$ws.Range("I57").Value = 1299
$ws.Range("J57").Value = 1300
$ws.Range("K57").Value = 3897
$ws.Range("L57").Value = 3900
$ws.Range("M57").Value = -3338
$ws.Range("N57").Value = -5018
$ws.Range("H94").Value = 8250
$ws.Range("J94").Value = 8250
$ws.Range("L94").Value = 24750
$ws.Range("N94").Value = -26102
$ws.Range("H104").Value = 6300.067
$ws.Range("I104").Value = 2071.5715
$ws.Range("J104").Value = 10000
$ws.Range("K104").Value = 6214.7145
$ws.Range("L104").Value = 30000
$ws.Range("M104").Value = -3593.7145
$ws.Range("N104").Value = -35242
$ws.Range("H109").Value = 48912.1
$ws.Range("I109").Value = 69375.14
$ws.Range("J109").Value = 1165
$ws.Range("K109").Value = 208125.42
$ws.Range("L109").Value = 3495
$ws.Range("M109").Value = -207085.42
$ws.Range("N109").Value = -5575
$ws.Range("H113").Value = 656.5833
$ws.Range("J113").Value = 734.8
$ws.Range("L113").Value = 2204.4
$ws.Range("N113").Value = -6544.4
$ws.Range("H121").Value = 243
$ws.Range("J121").Value = 500
$ws.Range("L121").Value = 1500
$ws.Range("N121").Value = -4120
$ws.Range("H131").Value = 1551.8
$ws.Range("I131").Value = 973
$ws.Range("K131").Value = 2919
$ws.Range("M131").Value = 2121

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4150
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4150
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 4150
$ws.Range("N80").Value = -6146
$ws.Range("H83").Value = 4150
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4150
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 20750
$ws.Range("N83").Value = -30734
$ws.Range("H94").Value = 24959
$ws.Range("J94").Value = 24959
$ws.Range("L94").Value = 24959
$ws.Range("N94").Value = -26311
$ws.Range("H102").Value = 1852.9
$ws.Range("I102").Value = 1991.25
$ws.Range("J102").Value = 1299.5
$ws.Range("K102").Value = 1991.25
$ws.Range("L102").Value = 1299.5
$ws.Range("M102").Value = -369.25
$ws.Range("N102").Value = -4543.5
$ws.Range("H122").Value = 1634.5333
$ws.Range("I122").Value = 1651.4
$ws.Range("J122").Value = 1600.8
$ws.Range("K122").Value = 4954.200000000001
$ws.Range("L122").Value = 4802.4
$ws.Range("M122").Value = -2504.200000000001
$ws.Range("N122").Value = -9702.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 18000
$ws.Range("I14").Value = 18000
$ws.Range("K14").Value = 18000
$ws.Range("M14").Value = -17828
$ws.Range("H82").Value = 4436.364
$ws.Range("I82").Value = 900
$ws.Range("K82").Value = 900
$ws.Range("M82").Value = -539
$ws.Range("H85").Value = 4436.364
$ws.Range("I85").Value = 900
$ws.Range("K85").Value = 900
$ws.Range("M85").Value = 348
$ws.Range("H136").Value = 4663.154
$ws.Range("I136").Value = 3311
$ws.Range("J136").Value = 4909
$ws.Range("K136").Value = 9933
$ws.Range("L136").Value = 14727
$ws.Range("M136").Value = -7383
$ws.Range("N136").Value = -19827

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 16549.857
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 28902.25
$ws.Range("K2").Value = 80
$ws.Range("L2").Value = 28902.25
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = -29126.25
$ws.Range("H41").Value = 22393.25
$ws.Range("J41").Value = 21498.75
$ws.Range("L41").Value = 21498.75
$ws.Range("N41").Value = -22278.75
$ws.Range("H107").Value = 2477
$ws.Range("I107").Value = 2723.1667
$ws.Range("K107").Value = 8169.500100000001
$ws.Range("M107").Value = -6249.500100000001
$ws.Range("H122").Value = 5850
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H130").Value = 7499.5
$ws.Range("J130").Value = 7499.5
$ws.Range("L130").Value = 7499.5
$ws.Range("N130").Value = -17539.5
$ws.Range("H132").Value = 1876.8572
$ws.Range("I132").Value = 1785
$ws.Range("K132").Value = 5355
$ws.Range("M132").Value = -2825
